$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 66; this shifts rows 66..94 down to 67..95.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the data of the "new" record.
# It mirrors what used to be row 65 (now shifted to row 66's neighbour logic),
# but carries its own (newer) price-reporting date.
$ws.Cells.Item(66, 1).Value = 11
$ws.Cells.Item(66, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(66, 3).Value = "Bíobío"
$ws.Cells.Item(66, 4).Value = 45119
$ws.Cells.Item(66, 5).Value = 8
$ws.Cells.Item(66, 6).Value = 100112013
$ws.Cells.Item(66, 7).Value = "Alcachofa"
$ws.Cells.Item(66, 8).Value = "Argentina(o)"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 100
$ws.Cells.Item(66, 11).Value = 14000
$ws.Cells.Item(66, 12).Value = 15000
$ws.Cells.Item(66, 13).Value = 14500
$ws.Cells.Item(66, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(66, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(66, 16).Value = 290
$ws.Cells.Item(66, 17).Value = 50
$ws.Cells.Item(66, 18).Value = "Hortaliza"
